$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Trf"
$ws.Cells.Item(2,3).Value = "Tfr2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.01309533333333333
$ws.Cells.Item(2,8).Value = 0.039286
$ws.Cells.Item(2,9).Value = 0.00009760639145116089
$ws.Cells.Item(2,10).Value = 0.00009760639145116087
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.010413
$ws.Cells.Item(2,14).Value = 0.031239
$ws.Cells.Item(2,15).Value = 0.01373154928834539
$ws.Cells.Item(2,16).Value = 0.01373154928834539
$ws.Cells.Item(2,17).Value = 0.000136361706
$ws.Cells.Item(2,18).Value = 0.001227255354
$ws.Cells.Item(2,19).Value = 0.00000134028697506915
$ws.Cells.Item(2,20).Value = 0.00000134028697506915

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Trf"
$ws.Cells.Item(3,3).Value = "Tfr2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.01309533333333333
$ws.Cells.Item(3,8).Value = 0.039286
$ws.Cells.Item(3,9).Value = 0.00009760639145116089
$ws.Cells.Item(3,10).Value = 0.00009760639145116087
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.7479136666666667
$ws.Cells.Item(3,14).Value = 2.243741
$ws.Cells.Item(3,15).Value = 0.9862684507116546
$ws.Cells.Item(3,16).Value = 0.9862684507116547
$ws.Cells.Item(3,17).Value = 0.009794178769555557
$ws.Cells.Item(3,18).Value = 0.088147608926
$ws.Cells.Item(3,19).Value = 0.00009626610447609173
$ws.Cells.Item(3,20).Value = 0.00009626610447609173

$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Trf"
$ws.Cells.Item(4,3).Value = "Tfr2"
$ws.Cells.Item(4,4).Value = "ECs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 132.5150096666667
$ws.Cells.Item(4,8).Value = 397.545029
$ws.Cells.Item(4,9).Value = 0.9877039077543427
$ws.Cells.Item(4,10).Value = 0.9877039077543425
$ws.Cells.Item(4,11).Value = 1
$ws.Cells.Item(4,12).Value = 0.3333333333333333
$ws.Cells.Item(4,13).Value = 0.010413
$ws.Cells.Item(4,14).Value = 0.031239
$ws.Cells.Item(4,15).Value = 0.01373154928834539
$ws.Cells.Item(4,16).Value = 0.01373154928834539
$ws.Cells.Item(4,17).Value = 1.379878795659
$ws.Cells.Item(4,18).Value = 12.418909160931
$ws.Cells.Item(4,19).Value = 0.01356270489162011
$ws.Cells.Item(4,20).Value = 0.01356270489162011

$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Trf"
$ws.Cells.Item(5,3).Value = "Tfr2"
$ws.Cells.Item(5,4).Value = "FAPs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 132.5150096666667
$ws.Cells.Item(5,8).Value = 397.545029
$ws.Cells.Item(5,9).Value = 0.9877039077543427
$ws.Cells.Item(5,10).Value = 0.9877039077543425
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.7479136666666667
$ws.Cells.Item(5,14).Value = 2.243741
$ws.Cells.Item(5,15).Value = 0.9862684507116546
$ws.Cells.Item(5,16).Value = 0.9862684507116547
$ws.Cells.Item(5,17).Value = 99.10978676816543
$ws.Cells.Item(5,18).Value = 891.9880809134889
$ws.Cells.Item(5,19).Value = 0.9741412028627225
$ws.Cells.Item(5,20).Value = 0.9741412028627224

$ws.Cells.Item(6,1).Value = "MuSCs"
$ws.Cells.Item(6,2).Value = "Trf"
$ws.Cells.Item(6,3).Value = "Tfr2"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 1.636606333333334
$ws.Cells.Item(6,8).Value = 4.909819000000001
$ws.Cells.Item(6,9).Value = 0.01219848585420627
$ws.Cells.Item(6,10).Value = 0.01219848585420626
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(6,12).Value = 0.3333333333333333
$ws.Cells.Item(6,13).Value = 0.010413
$ws.Cells.Item(6,14).Value = 0.031239
$ws.Cells.Item(6,15).Value = 0.01373154928834539
$ws.Cells.Item(6,16).Value = 0.01373154928834539
$ws.Cells.Item(6,17).Value = 0.017041981749
$ws.Cells.Item(6,18).Value = 0.153377835741
$ws.Cells.Item(6,19).Value = 0.0001675041097502174
$ws.Cells.Item(6,20).Value = 0.0001675041097502174

$ws.Cells.Item(7,1).Value = "MuSCs"
$ws.Cells.Item(7,2).Value = "Trf"
$ws.Cells.Item(7,3).Value = "Tfr2"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 1.636606333333334
$ws.Cells.Item(7,8).Value = 4.909819000000001
$ws.Cells.Item(7,9).Value = 0.01219848585420627
$ws.Cells.Item(7,10).Value = 0.01219848585420626
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.7479136666666667
$ws.Cells.Item(7,14).Value = 2.243741
$ws.Cells.Item(7,15).Value = 0.9862684507116546
$ws.Cells.Item(7,16).Value = 0.9862684507116547
$ws.Cells.Item(7,17).Value = 1.224040243653223
$ws.Cells.Item(7,18).Value = 11.016362192879
$ws.Cells.Item(7,19).Value = 0.01203098174445605
$ws.Cells.Item(7,20).Value = 0.01203098174445605
